$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 17.76613666666667
$ws.Range("H2").Value = 53.29841
$ws.Range("I2").Value = 0.7872390387208499
$ws.Range("J2").Value = 0.7872390387208499
$ws.Range("M2").Value = 49.89274333333334
$ws.Range("N2").Value = 149.67823
$ws.Range("O2").Value = 0.8663408689480834
$ws.Range("P2").Value = 0.8663408689480835
$ws.Range("Q2").Value = 886.4012967349223
$ws.Range("R2").Value = 7977.611670614301
$ws.Range("S2").Value = 0.682017352875275
$ws.Range("T2").Value = 0.682017352875275
$ws.Range("G3").Value = 17.76613666666667
$ws.Range("H3").Value = 53.29841
$ws.Range("I3").Value = 0.7872390387208499
$ws.Range("J3").Value = 0.7872390387208499
$ws.Range("O3").Value = 0.06984725491313053
$ws.Range("P3").Value = 0.06984725491313053
$ws.Range("Q3").Value = 71.46459268803555
$ws.Range("R3").Value = 643.18133419232
$ws.Range("S3").Value = 0.05498648581510304
$ws.Range("T3").Value = 0.05498648581510304
$ws.Range("G4").Value = 17.76613666666667
$ws.Range("H4").Value = 53.29841
$ws.Range("I4").Value = 0.7872390387208499
$ws.Range("J4").Value = 0.7872390387208499
$ws.Range("M4").Value = 1.266267666666667
$ws.Range("N4").Value = 3.798803
$ws.Range("O4").Value = 0.02198755485004457
$ws.Range("P4").Value = 0.02198755485004457
$ws.Range("Q4").Value = 22.49668442258112
$ws.Range("R4").Value = 202.47015980323
$ws.Range("S4").Value = 0.01730946154397104
$ws.Range("T4").Value = 0.01730946154397105
$ws.Range("G5").Value = 17.76613666666667
$ws.Range("H5").Value = 53.29841
$ws.Range("I5").Value = 0.7872390387208499
$ws.Range("J5").Value = 0.7872390387208499
$ws.Range("M5").Value = 0.2206823333333333
$ws.Range("N5").Value = 0.6620469999999999
$ws.Range("O5").Value = 0.003831942516052412
$ws.Range("P5").Value = 0.003831942516052413
$ws.Range("Q5").Value = 3.920672493918889
$ws.Range("R5").Value = 35.28605244527
$ws.Range("S5").Value = 0.003016654742770656
$ws.Range("T5").Value = 0.003016654742770656
$ws.Range("G6").Value = 17.76613666666667
$ws.Range("H6").Value = 53.29841
$ws.Range("I6").Value = 0.7872390387208499
$ws.Range("J6").Value = 0.7872390387208499
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.187988666666667
$ws.Range("N6").Value = 6.563966
$ws.Range("O6").Value = 0.03799237877268909
$ws.Range("P6").Value = 0.03799237877268909
$ws.Range("Q6").Value = 38.87210567711778
$ws.Range("R6").Value = 349.84895109406
$ws.Range("S6").Value = 0.02990908374373018
$ws.Range("T6").Value = 0.02990908374373018
$ws.Range("I7").Value = 0.03648413815195897
$ws.Range("J7").Value = 0.03648413815195897
$ws.Range("M7").Value = 49.89274333333334
$ws.Range("N7").Value = 149.67823
$ws.Range("O7").Value = 0.8663408689480834
$ws.Range("P7").Value = 0.8663408689480835
$ws.Range("Q7").Value = 41.07975567459111
$ws.Range("R7").Value = 369.71780107132
$ws.Range("S7").Value = 0.03160769994939006
$ws.Range("T7").Value = 0.03160769994939006
$ws.Range("I8").Value = 0.03648413815195897
$ws.Range("J8").Value = 0.03648413815195897
$ws.Range("O8").Value = 0.06984725491313053
$ws.Range("P8").Value = 0.06984725491313053
$ws.Range("S8").Value = 0.002548316897785749
$ws.Range("T8").Value = 0.002548316897785749
$ws.Range("I9").Value = 0.03648413815195897
$ws.Range("J9").Value = 0.03648413815195897
$ws.Range("M9").Value = 1.266267666666667
$ws.Range("N9").Value = 3.798803
$ws.Range("O9").Value = 0.02198755485004457
$ws.Range("P9").Value = 0.02198755485004457
$ws.Range("Q9").Value = 1.042595834383556
$ws.Range("R9").Value = 9.383362509452001
$ws.Range("S9").Value = 0.0008021969887728015
$ws.Range("T9").Value = 0.0008021969887728016
$ws.Range("I10").Value = 0.03648413815195897
$ws.Range("J10").Value = 0.03648413815195897
$ws.Range("M10").Value = 0.2206823333333333
$ws.Range("N10").Value = 0.6620469999999999
$ws.Range("O10").Value = 0.003831942516052412
$ws.Range("P10").Value = 0.003831942516052413
$ws.Range("Q10").Value = 0.1817013002164444
$ws.Range("R10").Value = 1.635311701948
$ws.Range("S10").Value = 0.0001398051201460215
$ws.Range("T10").Value = 0.0001398051201460215
$ws.Range("I11").Value = 0.03648413815195897
$ws.Range("J11").Value = 0.03648413815195897
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 2.187988666666667
$ws.Range("N11").Value = 6.563966
$ws.Range("O11").Value = 0.03799237877268909
$ws.Range("P11").Value = 0.03799237877268909
$ws.Range("Q11").Value = 1.801505265904889
$ws.Range("R11").Value = 16.213547393144
$ws.Range("S11").Value = 0.001386119195864342
$ws.Range("T11").Value = 0.001386119195864342
$ws.Range("G12").Value = 3.885299333333334
$ws.Range("H12").Value = 11.655898
$ws.Range("I12").Value = 0.1721623203571791
$ws.Range("J12").Value = 0.172162320357179
$ws.Range("M12").Value = 49.89274333333334
$ws.Range("N12").Value = 149.67823
$ws.Range("O12").Value = 0.8663408689480834
$ws.Range("P12").Value = 0.8663408689480835
$ws.Range("Q12").Value = 193.8482424111711
$ws.Range("R12").Value = 1744.63418170054
$ws.Range("S12").Value = 0.1491512542183568
$ws.Range("T12").Value = 0.1491512542183568
$ws.Range("G13").Value = 3.885299333333334
$ws.Range("H13").Value = 11.655898
$ws.Range("I13").Value = 0.1721623203571791
$ws.Range("J13").Value = 0.172162320357179
$ws.Range("O13").Value = 0.06984725491313053
$ws.Range("P13").Value = 0.06984725491313053
$ws.Range("Q13").Value = 15.62868391352178
$ws.Range("R13").Value = 140.658155221696
$ws.Range("S13").Value = 0.01202506547642393
$ws.Range("T13").Value = 0.01202506547642393
$ws.Range("G14").Value = 3.885299333333334
$ws.Range("H14").Value = 11.655898
$ws.Range("I14").Value = 0.1721623203571791
$ws.Range("J14").Value = 0.172162320357179
$ws.Range("M14").Value = 1.266267666666667
$ws.Range("N14").Value = 3.798803
$ws.Range("O14").Value = 0.02198755485004457
$ws.Range("P14").Value = 0.02198755485004457
$ws.Range("Q14").Value = 4.919828921121557
$ws.Range("R14").Value = 44.278460290094
$ws.Range("S14").Value = 0.003785428461964419
$ws.Range("T14").Value = 0.003785428461964419
$ws.Range("G15").Value = 3.885299333333334
$ws.Range("H15").Value = 11.655898
$ws.Range("I15").Value = 0.1721623203571791
$ws.Range("J15").Value = 0.172162320357179
$ws.Range("M15").Value = 0.2206823333333333
$ws.Range("N15").Value = 0.6620469999999999
$ws.Range("O15").Value = 0.003831942516052412
$ws.Range("P15").Value = 0.003831942516052413
$ws.Range("Q15").Value = 0.8574169225784445
$ws.Range("R15").Value = 7.716752303205999
$ws.Range("S15").Value = 0.0006597161150389102
$ws.Range("T15").Value = 0.0006597161150389102
$ws.Range("G16").Value = 3.885299333333334
$ws.Range("H16").Value = 11.655898
$ws.Range("I16").Value = 0.1721623203571791
$ws.Range("J16").Value = 0.172162320357179
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 2.187988666666667
$ws.Range("N16").Value = 6.563966
$ws.Range("O16").Value = 0.03799237877268909
$ws.Range("P16").Value = 0.03799237877268909
$ws.Range("Q16").Value = 8.500990907940889
$ws.Range("R16").Value = 76.508918171468
$ws.Range("S16").Value = 0.006540856085394988
$ws.Range("T16").Value = 0.006540856085394987
$ws.Range("G17").Value = 0.09285466666666665
$ws.Range("H17").Value = 0.278564
$ws.Range("I17").Value = 0.004114502770011991
$ws.Range("J17").Value = 0.004114502770011991
$ws.Range("M17").Value = 49.89274333333334
$ws.Range("N17").Value = 149.67823
$ws.Range("O17").Value = 0.8663408689480834
$ws.Range("P17").Value = 0.8663408689480835
$ws.Range("Q17").Value = 4.632774051302222
$ws.Range("R17").Value = 41.69496646172
$ws.Range("S17").Value = 0.003564561905061484
$ws.Range("T17").Value = 0.003564561905061484
$ws.Range("G18").Value = 0.09285466666666665
$ws.Range("H18").Value = 0.278564
$ws.Range("I18").Value = 0.004114502770011991
$ws.Range("J18").Value = 0.004114502770011991
$ws.Range("O18").Value = 0.06984725491313053
$ws.Range("P18").Value = 0.06984725491313053
$ws.Range("Q18").Value = 0.3735095061475555
$ws.Range("R18").Value = 3.361585555327999
$ws.Range("S18").Value = 0.0002873867238178092
$ws.Range("T18").Value = 0.0002873867238178092
$ws.Range("G19").Value = 0.09285466666666665
$ws.Range("H19").Value = 0.278564
$ws.Range("I19").Value = 0.004114502770011991
$ws.Range("J19").Value = 0.004114502770011991
$ws.Range("M19").Value = 1.266267666666667
$ws.Range("N19").Value = 3.798803
$ws.Range("O19").Value = 0.02198755485004457
$ws.Range("P19").Value = 0.02198755485004457
$ws.Range("Q19").Value = 0.1175788620991111
$ws.Range("R19").Value = 1.058209758892
$ws.Range("S19").Value = 0.00009046785533629894
$ws.Range("T19").Value = 0.00009046785533629895
$ws.Range("G20").Value = 0.09285466666666665
$ws.Range("H20").Value = 0.278564
$ws.Range("I20").Value = 0.004114502770011991
$ws.Range("J20").Value = 0.004114502770011991
$ws.Range("M20").Value = 0.2206823333333333
$ws.Range("N20").Value = 0.6620469999999999
$ws.Range("O20").Value = 0.003831942516052412
$ws.Range("P20").Value = 0.003831942516052413
$ws.Range("Q20").Value = 0.02049138450088888
$ws.Range("R20").Value = 0.184422460508
$ws.Range("S20").Value = 0.00001576653809682437
$ws.Range("T20").Value = 0.00001576653809682437
$ws.Range("G21").Value = 0.09285466666666665
$ws.Range("H21").Value = 0.278564
$ws.Range("I21").Value = 0.004114502770011991
$ws.Range("J21").Value = 0.004114502770011991
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 2.187988666666667
$ws.Range("N21").Value = 6.563966
$ws.Range("O21").Value = 0.03799237877268909
$ws.Range("P21").Value = 0.03799237877268909
$ws.Range("Q21").Value = 0.2031649583137777
$ws.Range("R21").Value = 1.828484624824
$ws.Range("S21").Value = 0.000156319747699574
$ws.Range("T21").Value = 0.000156319747699574
